$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44
$ws.Cells.Item(44,1).Value = 'EditDeal_TC001'
$ws.Cells.Item(44,1).VerticalAlignment = -4108
$ws.Cells.Item(44,2).Value = 'Validate whether a shipper user can edit deal in Opportunity.
a) Enter valid user id and Password and click Login(Shipper user).
b)Create deal and share with carrier user.
c) Check shared deal in opportunity
d) Edit deal and check updated details in opportunity
'
$ws.Cells.Item(44,2).WrapText = $true
$ws.Cells.Item(44,3).Value = 'NO'
$ws.Cells.Item(44,3).VerticalAlignment = -4108
$ws.Cells.Item(44,4).Value = 'Deal has been Edited'
$ws.Cells.Item(44,4).VerticalAlignment = -4108
$ws.Rows.Item(44).RowHeight = 105

# Row 45
$ws.Cells.Item(45,1).Value = 'EditDeal_TC002'
$ws.Cells.Item(45,1).VerticalAlignment = -4108
$ws.Cells.Item(45,2).Value = 'Validate whether a shipper admin can edit deal in Opportunity.
a) Enter valid user id and Password and click Login.
b)Create deal and share with carrier user.
c) Check shared deal in opportunity
d) Edit deal and check updated details in opportunity
'
$ws.Cells.Item(45,2).WrapText = $true
$ws.Cells.Item(45,3).Value = 'NO'
$ws.Cells.Item(45,3).VerticalAlignment = -4108
$ws.Cells.Item(45,4).Value = 'Deal has been Edited'
$ws.Cells.Item(45,4).VerticalAlignment = -4108
$ws.Rows.Item(45).RowHeight = 90

# Row 46
$ws.Cells.Item(46,1).Value = 'EditDeal_TC003'
$ws.Cells.Item(46,1).VerticalAlignment = -4108
$ws.Cells.Item(46,2).Value = 'Validate whether a carrier can edit deal in Opportunity.
a) Enter valid user id and Password and click Login.
b)Create deal and share with carrier user.
c) Check shared deal in opportunity
d) Edit deal and check updated details in opportunity
'
$ws.Cells.Item(46,2).WrapText = $true
$ws.Cells.Item(46,3).Value = 'NO'
$ws.Cells.Item(46,3).VerticalAlignment = -4108
$ws.Cells.Item(46,4).Value = 'Deal has been Edited'
$ws.Cells.Item(46,4).VerticalAlignment = -4108
$ws.Rows.Item(46).RowHeight = 90

# Row 47
$ws.Cells.Item(47,1).Value = 'ReShareDeal_TC001'
$ws.Cells.Item(47,1).VerticalAlignment = -4108
$ws.Cells.Item(47,2).Value = 'Validate whether a Shipper user can reshare deal in Opportunity.
a) Enter valid user id and Password and click Login.
b)Create deal and share with carrier user.
c) Check shared deal in opportunity
d) Again share deal with carrier user and check whether reshare popup displayed.
'
$ws.Cells.Item(47,2).WrapText = $true
$ws.Cells.Item(47,3).Value = 'NO'
$ws.Cells.Item(47,3).VerticalAlignment = -4108
$ws.Cells.Item(47,4).Value = 'Deal has been ReShared'
$ws.Cells.Item(47,4).VerticalAlignment = -4108
$ws.Rows.Item(47).RowHeight = 120

# Row 48
$ws.Cells.Item(48,1).Value = 'ReShareDeal_TC002'
$ws.Cells.Item(48,1).VerticalAlignment = -4108
$ws.Cells.Item(48,2).Value = 'Validate whether a carrier user can reshare deal in Opportunity.
a) Enter valid user id and Password and click Login.
b)Create deal and share with carrier user.
c) Check shared deal in opportunity
d) Again share deal with shipper user and check whether reshare popup displayed.
'
$ws.Cells.Item(48,2).WrapText = $true
$ws.Cells.Item(48,3).Value = 'NO'
$ws.Cells.Item(48,3).VerticalAlignment = -4108
$ws.Cells.Item(48,4).Value = 'Deal has been ReShared'
$ws.Cells.Item(48,4).VerticalAlignment = -4108
$ws.Rows.Item(48).RowHeight = 120

# Row 49
$ws.Cells.Item(49,1).Value = 'ReShareDeal_TC003'
$ws.Cells.Item(49,1).VerticalAlignment = -4108
$ws.Cells.Item(49,2).Value = 'Validate whether a Shipper admin can reshare deal in Opportunity.
a) Enter valid user id and Password and click Login.
b)Create deal and share with carrier user.
c) Check shared deal in opportunity
d) Again share deal with carrier user and check whether reshare popup displayed.
'
$ws.Cells.Item(49,2).WrapText = $true
$ws.Cells.Item(49,3).Value = 'NO'
$ws.Cells.Item(49,3).VerticalAlignment = -4108
$ws.Cells.Item(49,4).Value = 'Deal has been ReShared'
$ws.Cells.Item(49,4).VerticalAlignment = -4108
$ws.Rows.Item(49).RowHeight = 120

# Row 50
$ws.Cells.Item(50,1).Value = 'DiscardDeal_TC001'
$ws.Cells.Item(50,1).VerticalAlignment = -4108
$ws.Cells.Item(50,2).Value = 'Validate whehter Shipper user is able to discard Deal  on following conditions.
a) Launch application and login applcation as Shipper admin
b) Goto Deals and click on add icon
c) Set mandatory fileds and clik on next.
d) Click on Deals and Drafts.
e) Check whether the new deal widget added in Drafts.
f) Discard deal and check whether deal discarded successfully'
$ws.Cells.Item(50,2).WrapText = $true
$ws.Cells.Item(50,3).Value = 'NO'
$ws.Cells.Item(50,3).VerticalAlignment = -4108
$ws.Cells.Item(50,4).Value = 'Deal has been discarded'
$ws.Cells.Item(50,4).VerticalAlignment = -4108
$ws.Rows.Item(50).RowHeight = 120

# Row 51
$ws.Cells.Item(51,1).Value = 'DiscardDeal_TC002'
$ws.Cells.Item(51,1).VerticalAlignment = -4108
$ws.Cells.Item(51,2).Value = 'Validate whehter Shipper admin is able to discard Deal  on following conditions.
a) Launch application and login applcation as Shipper admin
b) Goto Deals and click on add icon
c) Set mandatory fileds and clik on next.
d) Click on Deals and Drafts.
e) Check whether the new deal widget added in Drafts.
f) Discard deal and check whether deal discarded successfully'
$ws.Cells.Item(51,2).WrapText = $true
$ws.Cells.Item(51,3).Value = 'Yes'
$ws.Cells.Item(51,3).VerticalAlignment = -4108
$ws.Cells.Item(51,4).Value = 'Deal has been discarded'
$ws.Cells.Item(51,4).VerticalAlignment = -4108
$ws.Rows.Item(51).RowHeight = 120

# Row 52
$ws.Cells.Item(52,1).Value = 'DiscardDeal_TC003'
$ws.Cells.Item(52,1).VerticalAlignment = -4108
$ws.Cells.Item(52,2).Value = 'Validate whehter carrier user is able to discard Deal  on following conditions.
a) Launch application and login applcation as Shipper admin
b) Goto Deals and click on add icon
c) Set mandatory fileds and clik on next.
d) Click on Deals and Drafts.
e) Check whether the new deal widget added in Drafts.
f) Discard deal and check whether deal discarded successfully'
$ws.Cells.Item(52,2).WrapText = $true
$ws.Cells.Item(52,3).Value = 'Yes'
$ws.Cells.Item(52,3).VerticalAlignment = -4108
$ws.Cells.Item(52,4).Value = 'Deal has been discarded'
$ws.Cells.Item(52,4).VerticalAlignment = -4108
$ws.Rows.Item(52).RowHeight = 120

# Update view state to match the selection/scroll position after the edit
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$win.ScrollColumn = 1
[void]$ws.Range("D44").Select()
